$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing columns right
$ws.Columns("A").Insert()

# Set the header for the newly inserted column
$ws.Range("A1").Value = "Stunden"
